$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New, re-sorted/deduplicated ticker list for column B (rows 2-59)
$tickers = @(
    "NSE:ADL",
    "NSE:ADVANIHOTR",
    "NSE:AIAENG",
    "NSE:APCL",
    "NSE:APTECHT",
    "NSE:ARTNIRMAN",
    "NSE:BODALCHEM",
    "NSE:CERA",
    "NSE:CHEMFAB",
    "NSE:CLEDUCATE",
    "NSE:CONCOR",
    "NSE:COSMOFIRST",
    "NSE:DCM",
    "NSE:DCW",
    "NSE:DEEPAKNTR",
    "NSE:DEN",
    "NSE:DIL",
    "NSE:DONEAR",
    "NSE:EMIL",
    "NSE:EPL",
    "NSE:ESSARSHPNG",
    "NSE:FINEORG",
    "NSE:GICHSGFIN",
    "NSE:GLOBUSSPR",
    "NSE:GMRINFRA",
    "NSE:GTLINFRA",
    "NSE:HATHWAY",
    "NSE:HCC",
    "NSE:HIKAL",
    "NSE:ICICIPRULI",
    "NSE:IFCI",
    "NSE:IGL",
    "NSE:IMFA",
    "NSE:JWL",
    "NSE:KAKATCEM",
    "NSE:KRITINUT",
    "NSE:LICHSGFIN",
    "NSE:LOWVOL",
    "NSE:MAHEPC",
    "NSE:MAXIND",
    "NSE:MMTC",
    "NSE:NATCOPHARM",
    "NSE:NATHBIOGEN",
    "NSE:NAVINFLUOR",
    "NSE:NDLVENTURE",
    "NSE:NECLIFE",
    "NSE:OMINFRAL",
    "NSE:ORIENTLTD",
    "NSE:PANSARI",
    "NSE:PDMJEPAPER",
    "NSE:PEL",
    "NSE:PNBHOUSING",
    "NSE:PRECWIRE",
    "NSE:PVRINOX",
    "NSE:RAILTEL",
    "NSE:RKDL",
    "NSE:RPOWER",
    "NSE:RUCHIRA"
)

for ($i = 0; $i -lt $tickers.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $tickers[$i]
}

# Column C updates
$ws.Cells.Item(2, 3).Value = "NSE:CAREERP"
$ws.Cells.Item(3, 3).Value = "NSE:KALYANIFRG"

# Column F updates
$ws.Cells.Item(2, 6).Value = "NSE:PEL"
$ws.Cells.Item(3, 6).Value = ""

# Delete rows 60-83 (old trailing rows that no longer exist)
$deleteRange = $ws.Range("A60:F83")
$deleteRange.EntireRow.Delete()
